# Applies the "Built site for gh-pages" edits to Carpinelli — Statement of Purpose.docx
#
# Strategy: use Range.Find/Execute for targeted text substitutions (phrase-level,
# not whole-paragraph rewrites, to stay as close as possible to the minimal set of
# wording changes), then fix up the two custom paragraph styles (delete the now
# unused "Abstract Title" style, and bump "Abstract" style's space-before).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# --- Date on the letterhead: "December 5, 2023" -> "December 14, 2023" ---
Replace-Text "December 5, 2023" "December 14, 2023"

# --- Drop the "(JSC)" parenthetical ---
Replace-Text "Space Center (JSC) Director." "Space Center Director."

# --- "As a flight dynamics engineer..." paragraph: phrase-level rewording ---
Replace-Text "As a flight dynamics engineer, I was primarily responsible for characterizing" `
             "As a flight dynamics engineer, I have characterized"
Replace-Text "and improving integrated spacecraft performance with numerical simulations. This work required all relevant dynamical effects to be modeled, including" `
             "and improved integrated spacecraft performance with numerical simulations. These simulations modeled all known dynamical effects, including"
Replace-Text "flexible structure & separation dynamics, propellant slosh, and sensor noise. I routinely modified dynamical models, and executed tens of thousands of Monte Carlo simulations to determine the vehicle performance’s sensitivity; one such simulation study found control parameter values which substantially improved vehicle performance. I independently" `
             "flexible structure & separation dynamics, propellant slosh, and sensor noise. I routinely executed tens of thousands of Monte Carlo simulations to determine the vehicle performance’s sensitivity to individually modified models. One such study found control parameter values which substantially improved vehicle performance; I independently"

# --- "In my three years at NASA..." paragraph ---
Replace-Text "In my three years at NASA, I have come to understand the space" `
             "After three years at NASA, I have come to understand the space"
Replace-Text "methods. My interest in space science has long been fueled by astronomers’ social media posts, and popular science magazines. I have been thrilled to recently learn that my technical skill-set may serve the cause of discovery through computational astronomy & astrophysics research." `
             "methods. Astronomers’ social media posts and popular science literature have long fueled my interested in space science, so I have been thrilled to learn that my technical skill-set can serve computational astronomy & astrophysics research."

# --- Punctuation tweaks ---
Replace-Text "; both" ": both"
Replace-Text "with my advisor, and quantify the impact of other feedback effects, e.g." `
             "with my advisor and quantify the impact of other feedback effects, e.g."
Replace-Text "including that within the" "including within the"
Replace-Text "exploring how Julia’s revolutionary modeling & simulation codes may" `
             "exploring how Julia’s revolutionary modeling and simulation codes may"
Replace-Text "accelerate astrophysical simulations, and aid projects such as Dr." `
             "accelerate astrophysical simulations and aid projects such as Dr."

# --- "With leaders in computational research..." closing paragraph rewrite ---
Replace-Text "With leaders in computational research, physics, mathematics, and other technical fields, MIT is uniquely suited to drive this moment in astronomy. I hope to have the opportunity to learn from this expertise. Thank you for your consideration." `
             "With leaders in computational research, physics, mathematics, and so many other technical fields, MIT is uniquely suited to drive this moment of interdiscplinary need in astronomy. I hope to have the opportunity to learn from this expertise as a graduate student in the MIT Department of Physics’ Doctoral Program. Thank you for your time, and your consideration."

# --- Styles: drop the unused "Abstract Title" style, and retune "Abstract" spacing ---
$abstractTitle = $d.Styles("Abstract Title")
$abstractTitle.Delete()

$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 15
